$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset to Normal style first so that newly-written cells don't inherit the
# sheet's column-level Text (@) number format; this matches the un-styled
# cells already used by columns B/C in the original sheet.
$ws.Range("D1:F7").Style = "Normal"
$ws.Range("A5:A7").Style = "Normal"

# ---- Row 1 (headers) ----
$ws.Range("D1").Value = " Oct 06"
$ws.Range("E1").Value = " Oct 07"
$ws.Range("F1").Value = " Oct 07"

# ---- Row 2 (Trucks) ----
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = 27
$ws.Range("F2").Value = 15

# ---- Row 3 (Dashboard) ----
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 5

# ---- Row 4 (Login) ----
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 2

# ---- Row 5 (Deleted Trailers) ----
$ws.Range("A5").Value = "Deleted Trailers"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0

# ---- Row 6 (Deleted Trucks) ----
$ws.Range("A6").Value = "Deleted Trucks"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0

# ---- Row 7 (Trailer) ----
$ws.Range("A7").Value = "Trailer"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
